# Q_device_1 : case 1 data update
# - widen columns A:B by one character (14.42578125 -> 15.42578125 stored width)
# - refresh the regression/residual values in rows 1-4
# - append a new row 5 of values, extending the used range to A1:B5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 14.67
$ws.Columns.Item(2).ColumnWidth = 14.67

$ws.Range("A1").Value = -0.045168227808977154
$ws.Range("B1").Value = -0.061518927093152299

$ws.Range("A2").Value = -0.0041884015802298334
$ws.Range("B2").Value = -0.0042166454971664918

$ws.Range("A3").Value = -0.024085272601921563
$ws.Range("B3").Value = -0.016707284668503995

$ws.Range("A4").Value = -0.0041258618718550836
$ws.Range("B4").Value = -0.0041402251454878899

$ws.Range("A5").Value = -0.059594559847822068
$ws.Range("B5").Value = -0.059639658827718024
